$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.788.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").Value = "'2.026.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = "'227.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.01%  '

$ws.Range("D6").Value = "'0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.10%  '

$ws.Range("D7").Value = "'60.14"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.55%  '

$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("D9").Value = "'0.381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").Value = "'0.0812"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.71%  '

$ws.Range("E11").Value = '  +0.72%  '

$ws.Range("D12").Value = "'14.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.44%  '

$ws.Range("D13").Value = "'2.326.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.28%  '

$ws.Range("D14").Value = "'21.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.64%  '

$ws.Range("D15").Value = "'0.753"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.53%  '

$ws.Range("D16").Value = "'5.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.35%  '

$ws.Range("D17").Value = "'2.048.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.85%  '

$ws.Range("D18").Value = "'37.730.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.91%  '

$ws.Range("D19").Value = "'6.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.44%  '

$ws.Range("D20").Value = "'69.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.82%  '

$ws.Range("D21").Value = "'0.0₃0825"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.10%  '

$ws.Range("D22").Value = "'223.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.18%  '

$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("E24").Value = '  -0.89%  '

$ws.Range("D25").Value = "'2.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.60%  '

$ws.Range("D26").Value = "'165.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.39%  '

$ws.Range("D27").Value = "'9.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("E28").Value = '  -3.31%  '

$ws.Range("D29").Value = "'18.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.48%  '

$ws.Range("D30").Value = "'1.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.58%  '

$ws.Range("E31").Value = '  +1.55%  '

$ws.Range("E32").Value = '  -0.91%  '

$ws.Range("D33").Value = "'2.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.71%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").Value = "'4.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.61%  '

$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = "'0.0602"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.53%  '

$ws.Range("D36").Value = "'6.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.76%  '

$ws.Range("D37").Value = "'2.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.00%  '

$ws.Range("D38").Value = "'3.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.12%  '

$ws.Range("E39").Value = '  -0.20%  '

$ws.Range("D40").Value = "'1.527.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.63%  '

$ws.Range("E41").Value = '  +0.81%  '

$ws.Range("D42").Value = "'96.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.97%  '

$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = "'16.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.82%  '

$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").Value = "'2.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.67%  '

$ws.Range("D45").Value = "'0.0918"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.29%  '

$ws.Range("D46").Value = "'1.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.33%  '

$ws.Range("D47").Value = "'4.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.47%  '

$ws.Range("D48").Value = "'2.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.54%  '

$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.57%  '

$ws.Range("E50").Value = '  -1.20%  '

$ws.Range("D51").Value = "'2.216.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.20%  '
